# Insert a new data row before row 125 (shifts existing rows 125-133 down
# to 126-134) and populate it with the new Berenjena price-report entry,
# matching the weekly update described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 125; this pushes rows 125-133 (and
# their formatting) down to 126-134 automatically.
$ws.Rows.Item(125).Insert()

# Populate the newly inserted row 125 with the new record.
$ws.Cells.Item(125, 1).Value = 7
$ws.Cells.Item(125, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(125, 3).Value = "Ñuble"
$ws.Cells.Item(125, 4).Value = 45267
$ws.Cells.Item(125, 5).Value = 16
$ws.Cells.Item(125, 6).Value = 100112001
$ws.Cells.Item(125, 7).Value = "Berenjena"
$ws.Cells.Item(125, 8).Value = "Sin especificar"
$ws.Cells.Item(125, 9).Value = "Primera"
$ws.Cells.Item(125, 10).Value = 60
$ws.Cells.Item(125, 11).Value = 10000
$ws.Cells.Item(125, 12).Value = 10000
$ws.Cells.Item(125, 13).Value = 10000
$ws.Cells.Item(125, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(125, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(125, 16).Value = 167
$ws.Cells.Item(125, 17).Value = 60
$ws.Cells.Item(125, 18).Value = "Hortaliza"
